# Swap the two worksheets' content:
#   - "hotel_info" (currently has the header row + one hotel data row, 9 cols)
#     becomes "review_info" (just the 25-column review_info header row).
#   - "review_info" (currently just a 25-column header row)
#     becomes "hotel_info" with a new "State" column inserted between
#     "Hotel_Name" and "City", populated with "Louisiana" for the one data row.
#
# Worksheet variables in this runtime resolve by *position*, not stable
# identity, so once we start moving sheets/cells around we always look the
# worksheets up again by their (still-original) names rather than relying on
# a previously captured reference.

$wb = $excel.ActiveWorkbook

$wsHotel  = $wb.Worksheets.Item("hotel_info")
$wsReview = $wb.Worksheets.Item("review_info")

# --- Stage the data that needs to move, using far-away scratch rows on the
#     destination sheet so we never overwrite something we still need to read.

# Stage review_info's header (A1:Y1) into hotel_info's scratch row 50.
$wsReview.Range("A1:Y1").Copy()
$wsHotel.Range("A50").PasteSpecial()

# Stage hotel_info's header+data into review_info's scratch rows 50:51,
# leaving a gap at column C for the new "State" column.
$wsHotel.Range("A1:B2").Copy()
$wsReview.Range("A50").PasteSpecial()
$wsHotel.Range("C1:I2").Copy()
$wsReview.Range("D50").PasteSpecial()
$wsReview.Range("C50").Value = "State"
$wsReview.Range("C51").Value = "Louisiana"

# --- Clear the old contents of both sheets.
$wsHotel.Range("A1:I2").Clear()
$wsReview.Range("A1:Y1").Clear()

# --- Move the staged data into its final A1 position on each sheet.
$wsHotel.Range("A50:Y50").Copy()
$wsHotel.Range("A1").PasteSpecial()
$wsHotel.Range("A50:Y50").Clear()

$wsReview.Range("A50:J51").Copy()
$wsReview.Range("A1").PasteSpecial()
$wsReview.Range("A50:J51").Clear()

# --- Rename the sheets to reflect their new contents. Swapping names needs
#     a temporary intermediate name so the two never collide.
$wsHotel.Name = "__tmp_swap__"
$wsReview.Name = "hotel_info"
$wsHotel.Name = "review_info"
